# Add a new "visibility_level" column (column Y) to the angels template
# worksheet, mirroring the existing header/data row pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (row 1) and sample data cell (row 2) for column Y.
$ws.Range("Y1").Value = "visibility_level"
$ws.Range("Y2").Value = "PRO"

# Give the new column an explicit width, matching the other data columns'
# custom-width styling (closest achievable value to 18.83203125 chars).
$ws.Range("Y1").EntireColumn.ColumnWidth = 18

$wb.Save()
